$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3672.3  # H76: 3795.2222 -> 3672.3
$ws.Cells.Item(76, 9).Value = 3340.625  # I76: 3455.8333 -> 3340.625
$ws.Cells.Item(76, 10).Value = 4999  # J76: 4474 -> 4999
$ws.Cells.Item(76, 11).Value = 3340.625  # K76: 3455.8333 -> 3340.625
$ws.Cells.Item(76, 12).Value = 4999  # L76: 4474 -> 4999
$ws.Cells.Item(76, 13).Value = -3025.625  # M76: -3140.8333 -> -3025.625
$ws.Cells.Item(76, 14).Value = -5629  # N76: -5104 -> -5629
$ws.Cells.Item(79, 8).Value = 3672.3  # H79: 3795.2222 -> 3672.3
$ws.Cells.Item(79, 9).Value = 3340.625  # I79: 3455.8333 -> 3340.625
$ws.Cells.Item(79, 10).Value = 4999  # J79: 4474 -> 4999
$ws.Cells.Item(79, 11).Value = 3340.625  # K79: 3455.8333 -> 3340.625
$ws.Cells.Item(79, 12).Value = 4999  # L79: 4474 -> 4999
$ws.Cells.Item(79, 13).Value = -2248.625  # M79: -2363.8333 -> -2248.625
$ws.Cells.Item(79, 14).Value = -7183  # N79: -6658 -> -7183
$ws.Cells.Item(87, 8).Value = 50000  # H87: 0 -> 50000
$ws.Cells.Item(87, 10).Value = 50000  # J87: 0 -> 50000
$ws.Cells.Item(87, 12).Value = 50000  # L87: 0 -> 50000
$ws.Cells.Item(87, 14).Value = -52496  # N87: None -> -52496
$ws.Cells.Item(90, 8).Value = 50000  # H90: 0 -> 50000
$ws.Cells.Item(90, 10).Value = 50000  # J90: 0 -> 50000
$ws.Cells.Item(90, 12).Value = 150000  # L90: 0 -> 150000
$ws.Cells.Item(90, 14).Value = -162480  # N90: None -> -162480
$ws.Cells.Item(92, 8).Value = 1726.2727  # H92: 1614.5385 -> 1726.2727
$ws.Cells.Item(92, 9).Value = 1600.8  # I92: 1400.125 -> 1600.8
$ws.Cells.Item(92, 10).Value = 1830.8334  # J92: 1957.6 -> 1830.8334
$ws.Cells.Item(92, 11).Value = 1600.8  # K92: 1400.125 -> 1600.8
$ws.Cells.Item(92, 12).Value = 1830.8334  # L92: 1957.6 -> 1830.8334
$ws.Cells.Item(92, 13).Value = -352.8  # M92: -152.125 -> -352.8
$ws.Cells.Item(92, 14).Value = -4326.8334  # N92: -4453.6 -> -4326.8334
$ws.Cells.Item(97, 8).Value = 811  # H97: 784 -> 811
$ws.Cells.Item(97, 10).Value = 811  # J97: 784 -> 811
$ws.Cells.Item(97, 12).Value = 2433  # L97: 2352 -> 2433
$ws.Cells.Item(97, 14).Value = -3425  # N97: -3344 -> -3425
$ws.Cells.Item(111, 8).Value = 908.1429000000001  # H111: 947.3077 -> 908.1429000000001
$ws.Cells.Item(111, 9).Value = 846.7273  # I111: 879.55554 -> 846.7273
$ws.Cells.Item(111, 10).Value = 1133.3334  # J111: 1099.75 -> 1133.3334
$ws.Cells.Item(111, 11).Value = 2540.1819  # K111: 2638.66662 -> 2540.1819
$ws.Cells.Item(111, 12).Value = 3400.0002  # L111: 3299.25 -> 3400.0002
$ws.Cells.Item(111, 13).Value = 526.8181  # M111: 428.33338 -> 526.8181
$ws.Cells.Item(111, 14).Value = -9534.0002  # N111: -9433.25 -> -9534.0002
$ws.Cells.Item(129, 8).Value = 2013.0769  # H129: 1924.7142 -> 2013.0769
$ws.Cells.Item(129, 9).Value = 1555.8334  # I129: 1495.8462 -> 1555.8334
$ws.Cells.Item(129, 11).Value = 4667.5002  # K129: 4487.5386 -> 4667.5002
$ws.Cells.Item(129, 13).Value = 332.4997999999996  # M129: 512.4614000000001 -> 332.4997999999996
$ws.Cells.Item(137, 8).Value = 471196.25  # H137: 406011.2 -> 471196.25
$ws.Cells.Item(137, 9).Value = 2206.9092  # I137: 1904.3572 -> 2206.9092
$ws.Cells.Item(137, 10).Value = 729140.4  # J137: 663170.0600000001 -> 729140.4
$ws.Cells.Item(137, 11).Value = 6620.7276  # K137: 5713.071599999999 -> 6620.7276
$ws.Cells.Item(137, 12).Value = 2187421.2  # L137: 1989510.18 -> 2187421.2
$ws.Cells.Item(137, 13).Value = -4070.7276  # M137: -3163.071599999999 -> -4070.7276
$ws.Cells.Item(137, 14).Value = -2192521.2  # N137: -1994610.18 -> -2192521.2

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 3903.875  # H5: 3661.6667 -> 3903.875
$ws.Cells.Item(5, 9).Value = 371.83334  # I5: 367.5 -> 371.83334
$ws.Cells.Item(5, 10).Value = 14500  # J5: 10250 -> 14500
$ws.Cells.Item(5, 11).Value = 371.83334  # K5: 367.5 -> 371.83334
$ws.Cells.Item(5, 12).Value = 14500  # L5: 10250 -> 14500
$ws.Cells.Item(5, 13).Value = -259.83334  # M5: -255.5 -> -259.83334
$ws.Cells.Item(5, 14).Value = -14724  # N5: -10474 -> -14724
$ws.Cells.Item(45, 8).Value = 9622569  # H45: 10424233 -> 9622569
$ws.Cells.Item(45, 9).Value = 2625  # I45: 2633.6667 -> 2625
$ws.Cells.Item(45, 11).Value = 2625  # K45: 2633.6667 -> 2625
$ws.Cells.Item(45, 13).Value = -2248  # M45: -2256.6667 -> -2248
$ws.Cells.Item(97, 8).Value = 0  # H97: 1272.5 -> 0
$ws.Cells.Item(97, 9).Value = 0  # I97: 1250.1 -> 0
$ws.Cells.Item(97, 10).Value = 0  # J97: 1384.5 -> 0
$ws.Cells.Item(97, 11).Value = 0  # K97: 1250.1 -> 0
$ws.Cells.Item(97, 12).Value = 0  # L97: 1384.5 -> 0
$ws.Cells.Item(97, 13).ClearContents()  # M97: was -754.0999999999999
$ws.Cells.Item(97, 14).ClearContents()  # N97: was -2376.5
$ws.Cells.Item(122, 8).Value = 1622.4736  # H122: 1748.9375 -> 1622.4736
$ws.Cells.Item(122, 9).Value = 1380.6364  # I122: 1537.5555 -> 1380.6364
$ws.Cells.Item(122, 10).Value = 1955  # J122: 2020.7142 -> 1955
$ws.Cells.Item(122, 11).Value = 4141.9092  # K122: 4612.666499999999 -> 4141.9092
$ws.Cells.Item(122, 12).Value = 5865  # L122: 6062.142599999999 -> 5865
$ws.Cells.Item(122, 13).Value = -1691.9092  # M122: -2162.666499999999 -> -1691.9092
$ws.Cells.Item(122, 14).Value = -10765  # N122: -10962.1426 -> -10765

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 3903.875  # H4: 3661.6667 -> 3903.875
$ws.Cells.Item(4, 9).Value = 371.83334  # I4: 367.5 -> 371.83334
$ws.Cells.Item(4, 10).Value = 14500  # J4: 10250 -> 14500
$ws.Cells.Item(4, 11).Value = 371.83334  # K4: 367.5 -> 371.83334
$ws.Cells.Item(4, 12).Value = 14500  # L4: 10250 -> 14500
$ws.Cells.Item(4, 13).Value = -256.83334  # M4: -252.5 -> -256.83334
$ws.Cells.Item(4, 14).Value = -14730  # N4: -10480 -> -14730

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1548.5555  # H31: 1902.2 -> 1548.5555
$ws.Cells.Item(31, 9).Value = 1548.5555  # I31: 1627.75 -> 1548.5555
$ws.Cells.Item(31, 10).Value = 0  # J31: 3000 -> 0
$ws.Cells.Item(31, 11).Value = 1548.5555  # K31: 1627.75 -> 1548.5555
$ws.Cells.Item(31, 12).Value = 0  # L31: 3000 -> 0
$ws.Cells.Item(31, 13).Value = -1253.5555  # M31: -1332.75 -> -1253.5555
$ws.Cells.Item(31, 14).ClearContents()  # N31: was -3590
$ws.Cells.Item(34, 8).Value = 1548.5555  # H34: 1902.2 -> 1548.5555
$ws.Cells.Item(34, 9).Value = 1548.5555  # I34: 1627.75 -> 1548.5555
$ws.Cells.Item(34, 10).Value = 0  # J34: 3000 -> 0
$ws.Cells.Item(34, 11).Value = 1548.5555  # K34: 1627.75 -> 1548.5555
$ws.Cells.Item(34, 12).Value = 0  # L34: 3000 -> 0
$ws.Cells.Item(34, 13).Value = -1346.5555  # M34: -1425.75 -> -1346.5555
$ws.Cells.Item(34, 14).ClearContents()  # N34: was -3404
$ws.Cells.Item(93, 8).Value = 34400  # H93: 30775 -> 34400
$ws.Cells.Item(93, 9).Value = 34400  # I93: 30775 -> 34400
$ws.Cells.Item(93, 11).Value = 34400  # K93: 30775 -> 34400
$ws.Cells.Item(93, 13).Value = -32528  # M93: -28903 -> -32528
$ws.Cells.Item(141, 8).Value = 100000  # H141: 275000 -> 100000
$ws.Cells.Item(141, 10).Value = 100000  # J141: 275000 -> 100000
$ws.Cells.Item(141, 12).Value = 100000  # L141: 275000 -> 100000
$ws.Cells.Item(141, 14).Value = -110360  # N141: -285360 -> -110360

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 5613.6924  # H12: 5462.7144 -> 5613.6924
$ws.Cells.Item(12, 9).Value = 5613.6924  # I12: 5462.7144 -> 5613.6924
$ws.Cells.Item(12, 11).Value = 5613.6924  # K12: 5462.7144 -> 5613.6924
$ws.Cells.Item(12, 13).Value = -5473.6924  # M12: -5322.7144 -> -5473.6924
$ws.Cells.Item(122, 8).Value = 6750.1333  # H122: 7521.3335 -> 6750.1333
$ws.Cells.Item(122, 9).Value = 7532.125  # I122: 8179.5713 -> 7532.125
$ws.Cells.Item(122, 10).Value = 5856.4287  # J122: 6599.8 -> 5856.4287
$ws.Cells.Item(122, 11).Value = 22596.375  # K122: 24538.7139 -> 22596.375
$ws.Cells.Item(122, 12).Value = 17569.2861  # L122: 19799.4 -> 17569.2861
$ws.Cells.Item(122, 13).Value = -20146.375  # M122: -22088.7139 -> -20146.375
$ws.Cells.Item(122, 14).Value = -22469.2861  # N122: -24699.4 -> -22469.2861

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3510.8  # H7: 3675.7778 -> 3510.8
$ws.Cells.Item(7, 9).Value = 2814.7  # I7: 3011.875 -> 2814.7
$ws.Cells.Item(7, 11).Value = 2814.7  # K7: 3011.875 -> 2814.7
$ws.Cells.Item(7, 13).Value = -2702.7  # M7: -2899.875 -> -2702.7
$ws.Cells.Item(16, 8).Value = 2627.15  # H16: 2632.35 -> 2627.15
$ws.Cells.Item(16, 9).Value = 2277.7334  # I16: 2284.6667 -> 2277.7334
$ws.Cells.Item(16, 11).Value = 2277.7334  # K16: 2284.6667 -> 2277.7334
$ws.Cells.Item(16, 13).Value = -2107.7334  # M16: -2114.6667 -> -2107.7334
$ws.Cells.Item(20, 8).Value = 13857.857  # H20: 13375.5 -> 13857.857
$ws.Cells.Item(20, 10).Value = 14500  # J20: 13857 -> 14500
$ws.Cells.Item(20, 12).Value = 14500  # L20: 13857 -> 14500
$ws.Cells.Item(20, 14).Value = -14952  # N20: -14309 -> -14952
$ws.Cells.Item(55, 8).Value = 2644.4348  # H55: 2849.9048 -> 2644.4348
$ws.Cells.Item(55, 9).Value = 1897.9166  # I55: 2180.1 -> 1897.9166
$ws.Cells.Item(55, 11).Value = 1897.9166  # K55: 2180.1 -> 1897.9166
$ws.Cells.Item(55, 13).Value = -1724.9166  # M55: -2007.1 -> -1724.9166
$ws.Cells.Item(61, 8).Value = 4111.154  # H61: 4324.5835 -> 4111.154
$ws.Cells.Item(61, 9).Value = 4472.778  # I61: 4838.125 -> 4472.778
$ws.Cells.Item(61, 11).Value = 4472.778  # K61: 4838.125 -> 4472.778
$ws.Cells.Item(61, 13).Value = -4270.778  # M61: -4636.125 -> -4270.778
$ws.Cells.Item(82, 8).Value = 2382.3044  # H82: 2251.72 -> 2382.3044
$ws.Cells.Item(82, 9).Value = 2279.3333  # I82: 2099.4119 -> 2279.3333
$ws.Cells.Item(82, 11).Value = 2279.3333  # K82: 2099.4119 -> 2279.3333
$ws.Cells.Item(82, 13).Value = -1918.3333  # M82: -1738.4119 -> -1918.3333
$ws.Cells.Item(85, 8).Value = 2382.3044  # H85: 2251.72 -> 2382.3044
$ws.Cells.Item(85, 9).Value = 2279.3333  # I85: 2099.4119 -> 2279.3333
$ws.Cells.Item(85, 11).Value = 2279.3333  # K85: 2099.4119 -> 2279.3333
$ws.Cells.Item(85, 13).Value = -1031.3333  # M85: -851.4119000000001 -> -1031.3333
$ws.Cells.Item(113, 8).Value = 4111.154  # H113: 4324.5835 -> 4111.154
$ws.Cells.Item(113, 9).Value = 4472.778  # I113: 4838.125 -> 4472.778
$ws.Cells.Item(113, 11).Value = 4472.778  # K113: 4838.125 -> 4472.778
$ws.Cells.Item(113, 13).Value = -2302.778  # M113: -2668.125 -> -2302.778
$ws.Cells.Item(126, 8).Value = 3510.8  # H126: 3675.7778 -> 3510.8
$ws.Cells.Item(126, 9).Value = 2814.7  # I126: 3011.875 -> 2814.7
$ws.Cells.Item(126, 11).Value = 8444.099999999999  # K126: 9035.625 -> 8444.099999999999
$ws.Cells.Item(126, 13).Value = -5974.099999999999  # M126: -6565.625 -> -5974.099999999999
$ws.Cells.Item(136, 8).Value = 4181.8823  # H136: 3648.0476 -> 4181.8823
$ws.Cells.Item(136, 9).Value = 4224.625  # I136: 3655.55 -> 4224.625
$ws.Cells.Item(136, 11).Value = 12673.875  # K136: 10966.65 -> 12673.875
$ws.Cells.Item(136, 13).Value = -10123.875  # M136: -8416.650000000001 -> -10123.875

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 21399.6  # H15: 24999.5 -> 21399.6
$ws.Cells.Item(15, 10).Value = 21399.6  # J15: 24999.5 -> 21399.6
$ws.Cells.Item(15, 12).Value = 21399.6  # L15: 24999.5 -> 21399.6
$ws.Cells.Item(15, 14).Value = -21975.6  # N15: -25575.5 -> -21975.6
$ws.Cells.Item(46, 8).Value = 92396.5  # H46: 92396.664 -> 92396.5
$ws.Cells.Item(46, 10).Value = 92396.5  # J46: 92396.664 -> 92396.5
$ws.Cells.Item(46, 12).Value = 92396.5  # L46: 92396.664 -> 92396.5
$ws.Cells.Item(46, 14).Value = -92858.5  # N46: -92858.664 -> -92858.5
$ws.Cells.Item(112, 8).Value = 43499  # H112: 43999 -> 43499
$ws.Cells.Item(112, 10).Value = 43499  # J112: 43999 -> 43499
$ws.Cells.Item(112, 12).Value = 43499  # L112: 43999 -> 43499
$ws.Cells.Item(112, 14).Value = -46453  # N112: -46953 -> -46453
$ws.Cells.Item(119, 8).Value = 45700  # H119: 47449.75 -> 45700
$ws.Cells.Item(119, 10).Value = 45700  # J119: 47449.75 -> 45700
$ws.Cells.Item(119, 12).Value = 45700  # L119: 47449.75 -> 45700
$ws.Cells.Item(119, 14).Value = -55376  # N119: -57125.75 -> -55376
$ws.Cells.Item(129, 8).Value = 69998  # H129: 70000 -> 69998
$ws.Cells.Item(129, 10).Value = 69998  # J129: 70000 -> 69998
$ws.Cells.Item(129, 12).Value = 69998  # L129: 70000 -> 69998
$ws.Cells.Item(129, 14).Value = -79998  # N129: -80000 -> -79998
$ws.Cells.Item(134, 8).Value = 92396.5  # H134: 92396.664 -> 92396.5
$ws.Cells.Item(134, 10).Value = 92396.5  # J134: 92396.664 -> 92396.5
$ws.Cells.Item(134, 12).Value = 277189.5  # L134: 277189.992 -> 277189.5
$ws.Cells.Item(134, 14).Value = -282259.5  # N134: -282259.992 -> -282259.5
